# "Changes of 29th march 2022"
#
# ShipmentTracking numbers in P2/P3 were re-run and got new FedEx tracking
# numbers. The source cells are General-formatted but hold their tracking
# numbers as literal text (shared strings) rather than numbers, and that
# must be preserved (no NumberFormat/style change on P2/P3).
#
# A plain `.Value = "320018151874"` gets auto-coerced to a Number by this
# engine because the string is all-digits, which would change the cell's
# XML shape (t="s" -> numeric <v>) and isn't what the diff shows. Using a
# helper cell with a text formula and pasting its *value* back in keeps the
# destination cell's existing format/style untouched while still writing a
# genuine shared-string (t="s") literal, matching the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cell, [string]$text) {
    $helper = $ws.Range("ZZ1")
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $helper.Value = $null
}

Set-TextValue $ws.Range("P2") "320018151874"
Set-TextValue $ws.Range("P3") "320018151885"

$excel.CutCopyMode = $false
